# ----------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at the
#    top of the data block (existing 2022-Q3 / Q2 / Q1 rows shift down).
# 2. Insert a brand-new "2022-Q4" worksheet right after "总计" (so the sheet
#    order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1) and populate it
#    with the quarter's fund-holding detail rows.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# --- 1. "总计" sheet: insert the new 2022-Q4 row ---------------------------
$total.Rows.Item(2).Insert()

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# The freshly inserted row 2 inherited row 1's (header) bold/border style and
# row 3 lost its index-column style during the shift; re-sync both rows with
# row 3's original untouched formatting (A col = index style, B:D = plain).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 21
$total.Range("D2").Value = 2.79

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 19
$total.Range("D3").Value = 1.7

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.33

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.55

# --- 2. Insert the new "2022-Q4" worksheet right after "总计" --------------
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q4"

# Header row (copy the bold/bordered header style from the "总计" sheet so no
# brand-new style entries are minted).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$total.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "290011"
$newSheet.Range("C2").Value = "泰信中小盘精选混合"
$newSheet.Range("D2").Value = "14.75"
$newSheet.Range("E2").Value = "94.33"
$newSheet.Range("F2").Value = "10.07"
$newSheet.Range("G2").Value = "1.4853"
$newSheet.Range("B2:G2").Style = "Normal"
$newSheet.Range("H2").Value = 1
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "002580"
$newSheet.Range("C3").Value = "泰信鑫选灵活配置混合C"
$newSheet.Range("D3").Value = "1.62"
$newSheet.Range("E3").Value = "93.92"
$newSheet.Range("F3").Value = "9.98"
$newSheet.Range("G3").Value = "0.1617"
$newSheet.Range("B3:G3").Style = "Normal"
$newSheet.Range("H3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4:G4").NumberFormat = "@"
$newSheet.Range("B4").Value = "000742"
$newSheet.Range("C4").Value = "国泰新经济灵活配置混合A"
$newSheet.Range("D4").Value = "5.03"
$newSheet.Range("E4").Value = "86.40"
$newSheet.Range("F4").Value = "3.12"
$newSheet.Range("G4").Value = "0.1569"
$newSheet.Range("B4:G4").Style = "Normal"
$newSheet.Range("H4").Value = 9
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5:G5").NumberFormat = "@"
$newSheet.Range("B5").Value = "006058"
$newSheet.Range("C5").Value = "民生加银新兴成长混合"
$newSheet.Range("D5").Value = "3.44"
$newSheet.Range("E5").Value = "86.91"
$newSheet.Range("F5").Value = "3.91"
$newSheet.Range("G5").Value = "0.1345"
$newSheet.Range("B5:G5").Style = "Normal"
$newSheet.Range("H5").Value = 7
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6:G6").NumberFormat = "@"
$newSheet.Range("B6").Value = "009234"
$newSheet.Range("C6").Value = "鹏华优质企业混合"
$newSheet.Range("D6").Value = "3.21"
$newSheet.Range("E6").Value = "80.54"
$newSheet.Range("F6").Value = "3.55"
$newSheet.Range("G6").Value = "0.1140"
$newSheet.Range("B6:G6").Style = "Normal"
$newSheet.Range("H6").Value = 8
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7:G7").NumberFormat = "@"
$newSheet.Range("B7").Value = "008811"
$newSheet.Range("C7").Value = "鹏华科技创新混合"
$newSheet.Range("D7").Value = "2.83"
$newSheet.Range("E7").Value = "89.03"
$newSheet.Range("F7").Value = "3.54"
$newSheet.Range("G7").Value = "0.1002"
$newSheet.Range("B7:G7").Style = "Normal"
$newSheet.Range("H7").Value = 5
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8:G8").NumberFormat = "@"
$newSheet.Range("B8").Value = "001970"
$newSheet.Range("C8").Value = "泰信鑫选灵活配置混合A"
$newSheet.Range("D8").Value = "0.94"
$newSheet.Range("E8").Value = "93.92"
$newSheet.Range("F8").Value = "9.98"
$newSheet.Range("G8").Value = "0.0938"
$newSheet.Range("B8:G8").Style = "Normal"
$newSheet.Range("H8").Value = 1
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9:G9").NumberFormat = "@"
$newSheet.Range("B9").Value = "501200"
$newSheet.Range("C9").Value = "民生加银科技创新 3 年封闭混合"
$newSheet.Range("D9").Value = "2.45"
$newSheet.Range("E9").Value = "86.66"
$newSheet.Range("F9").Value = "3.79"
$newSheet.Range("G9").Value = "0.0929"
$newSheet.Range("B9:G9").Style = "Normal"
$newSheet.Range("H9").Value = 9
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10:G10").NumberFormat = "@"
$newSheet.Range("B10").Value = "005819"
$newSheet.Range("C10").Value = "国泰优势行业混合A"
$newSheet.Range("D10").Value = "2.24"
$newSheet.Range("E10").Value = "90.14"
$newSheet.Range("F10").Value = "3.33"
$newSheet.Range("G10").Value = "0.0746"
$newSheet.Range("B10:G10").Style = "Normal"
$newSheet.Range("H10").Value = 8
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11:G11").NumberFormat = "@"
$newSheet.Range("B11").Value = "010912"
$newSheet.Range("C11").Value = "国泰成长价值混合A"
$newSheet.Range("D11").Value = "2.13"
$newSheet.Range("E11").Value = "86.46"
$newSheet.Range("F11").Value = "3.14"
$newSheet.Range("G11").Value = "0.0669"
$newSheet.Range("B11:G11").Style = "Normal"
$newSheet.Range("H11").Value = 9
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12:G12").NumberFormat = "@"
$newSheet.Range("B12").Value = "014606"
$newSheet.Range("C12").Value = "招商高端装备混合A"
$newSheet.Range("D12").Value = "2.02"
$newSheet.Range("E12").Value = "94.44"
$newSheet.Range("F12").Value = "3.12"
$newSheet.Range("G12").Value = "0.0630"
$newSheet.Range("B12:G12").Style = "Normal"
$newSheet.Range("H12").Value = 10
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13:G13").NumberFormat = "@"
$newSheet.Range("B13").Value = "014686"
$newSheet.Range("C13").Value = "招商核心装备混合A"
$newSheet.Range("D13").Value = "1.90"
$newSheet.Range("E13").Value = "91.32"
$newSheet.Range("F13").Value = "2.89"
$newSheet.Range("G13").Value = "0.0549"
$newSheet.Range("B13:G13").Style = "Normal"
$newSheet.Range("H13").Value = 10
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14:G14").NumberFormat = "@"
$newSheet.Range("B14").Value = "011712"
$newSheet.Range("C14").Value = "大摩万众创新混合C"
$newSheet.Range("D14").Value = "0.71"
$newSheet.Range("E14").Value = "93.83"
$newSheet.Range("F14").Value = "7.23"
$newSheet.Range("G14").Value = "0.0513"
$newSheet.Range("B14:G14").Style = "Normal"
$newSheet.Range("H14").Value = 4
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15:G15").NumberFormat = "@"
$newSheet.Range("B15").Value = "014607"
$newSheet.Range("C15").Value = "招商高端装备混合C"
$newSheet.Range("D15").Value = "1.62"
$newSheet.Range("E15").Value = "94.44"
$newSheet.Range("F15").Value = "3.12"
$newSheet.Range("G15").Value = "0.0505"
$newSheet.Range("B15:G15").Style = "Normal"
$newSheet.Range("H15").Value = 10
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16:G16").NumberFormat = "@"
$newSheet.Range("B16").Value = "002885"
$newSheet.Range("C16").Value = "大摩万众创新混合A"
$newSheet.Range("D16").Value = "0.36"
$newSheet.Range("E16").Value = "93.83"
$newSheet.Range("F16").Value = "7.23"
$newSheet.Range("G16").Value = "0.0260"
$newSheet.Range("B16:G16").Style = "Normal"
$newSheet.Range("H16").Value = 4
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17:G17").NumberFormat = "@"
$newSheet.Range("B17").Value = "014687"
$newSheet.Range("C17").Value = "招商核心装备混合C"
$newSheet.Range("D17").Value = "0.84"
$newSheet.Range("E17").Value = "91.32"
$newSheet.Range("F17").Value = "2.89"
$newSheet.Range("G17").Value = "0.0243"
$newSheet.Range("B17:G17").Style = "Normal"
$newSheet.Range("H17").Value = 10
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18:G18").NumberFormat = "@"
$newSheet.Range("B18").Value = "015585"
$newSheet.Range("C18").Value = "国泰优势行业混合C"
$newSheet.Range("D18").Value = "0.70"
$newSheet.Range("E18").Value = "90.14"
$newSheet.Range("F18").Value = "3.33"
$newSheet.Range("G18").Value = "0.0233"
$newSheet.Range("B18:G18").Style = "Normal"
$newSheet.Range("H18").Value = 8
$newSheet.Range("A19").Value = 17
$newSheet.Range("B19:G19").NumberFormat = "@"
$newSheet.Range("B19").Value = "006072"
$newSheet.Range("C19").Value = "民生加银创新成长混合A"
$newSheet.Range("D19").Value = "0.40"
$newSheet.Range("E19").Value = "91.73"
$newSheet.Range("F19").Value = "3.08"
$newSheet.Range("G19").Value = "0.0123"
$newSheet.Range("B19:G19").Style = "Normal"
$newSheet.Range("H19").Value = 10
$newSheet.Range("A20").Value = 18
$newSheet.Range("B20:G20").NumberFormat = "@"
$newSheet.Range("B20").Value = "014989"
$newSheet.Range("C20").Value = "国泰新经济灵活配置混合C"
$newSheet.Range("D20").Value = "0.14"
$newSheet.Range("E20").Value = "86.40"
$newSheet.Range("F20").Value = "3.12"
$newSheet.Range("G20").Value = "0.0044"
$newSheet.Range("B20:G20").Style = "Normal"
$newSheet.Range("H20").Value = 9
$newSheet.Range("A21").Value = 19
$newSheet.Range("B21:G21").NumberFormat = "@"
$newSheet.Range("B21").Value = "010913"
$newSheet.Range("C21").Value = "国泰成长价值混合C"
$newSheet.Range("D21").Value = "0.11"
$newSheet.Range("E21").Value = "86.46"
$newSheet.Range("F21").Value = "3.14"
$newSheet.Range("G21").Value = "0.0035"
$newSheet.Range("B21:G21").Style = "Normal"
$newSheet.Range("H21").Value = 9
$newSheet.Range("A22").Value = 20
$newSheet.Range("B22:G22").NumberFormat = "@"
$newSheet.Range("B22").Value = "014929"
$newSheet.Range("C22").Value = "民生加银创新成长混合C"
$newSheet.Range("D22").Value = "0.01"
$newSheet.Range("E22").Value = "91.73"
$newSheet.Range("F22").Value = "3.08"
$newSheet.Range("G22").Value = "0.0003"
$newSheet.Range("B22:G22").Style = "Normal"
$newSheet.Range("H22").Value = 10

# Apply the A-column "index" style (same style used by the 总计 sheet's A
# column / the other quarter sheets) across all 21 data rows in one shot.
$total.Range("A2").Copy()
$newSheet.Range("A2:A22").PasteSpecial(-4122)

Write-Output "2022-Q4 sheet + summary row added"
